# "shop on main scene working"
# Adds a "Button Shopping Cart" attribution row to the Pistures asset-source
# sheet, inserting a couple of spacer rows above it (matching the row
# layout produced by Excel when a new block of rows is inserted in the
# middle of the existing attribution table), and moves the active
# selection to A7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the "Button Music" row (old row 6), pushing
# the "Button Music" / "Button Setting" rows down from 6-7 to 8-9.
$ws.Rows("6:7").Insert()

# Insert one more row after the (now shifted) "Button Setting" row (row 9)
# to make room for the new "Button Shopping Cart" entry at row 10. This
# also pushes the remaining attribution rows further down (old 9-11 -> 12-14,
# old 13 -> 16), keeping the same blank-row gaps that separated sections.
$ws.Rows("10:10").Insert()

# Populate the new "Button Shopping Cart" attribution row.
$ws.Range("A10").Value = "Button Shopping Cart"
$ws.Range("C10").Value = "Attribution, Commercially, Modify"
$ws.Range("D10").Value = "CC BY 3.0  "
$ws.Range("E10").Value = "http://creativecommons.org/licenses/by/3.0/"
$ws.Range("F10").Value = "http://www.flaticon.com/authors/google"

# Match style of the new spacer rows' C column (style index used by the
# rest of the attribution "Type" column) -- inherited automatically by the
# row insert above, so nothing else required there.

# Move the active selection to A7 (inside the new spacer block).
$ws.Range("A7").Select()
